$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Company")
$ws.Range("A2").Value = "CapProviderTestCompany"
